$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.305.44'
$ws.Range("E2").Value = '  -0.17%  '

$ws.Range("D3").Value = '2.648.77'
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.28'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.23'
$ws.Range("E6").Value = '  +1.74%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("E9").Value = '  +2.85%  '

$ws.Range("E10").Value = '  -1.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.25'
$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("E12").Value = '  +0.80%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.05'
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000189'
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("D15").Value = '3.130.23'
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("D16").Value = '68.199.44'
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("D17").Value = '2.655.96'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.39'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '364.49'
$ws.Range("E19").Value = '  -0.42%  '

$ws.Range("E20").Value = '  -0.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.43'
$ws.Range("E21").Value = '  +4.06%  '

$ws.Range("E22").Value = '  -1.07%  '

$ws.Range("E23").Value = '  -2.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.29'
$ws.Range("E24").Value = '  +2.60%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.75'
$ws.Range("E26").Value = '  -2.59%  '

$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("E28").Value = '  +0.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '559.72'
$ws.Range("E30").Value = '  -2.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.06'
$ws.Range("E31").Value = '  +0.83%  '

$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("E33").Value = '  +0.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.129'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("E36").Value = '  +0.51%  '

$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.85'
$ws.Range("E37").Value = '  +3.35%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.52'
$ws.Range("E38").Value = '  +0.28%  '

$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.88'
$ws.Range("E40").Value = '  -2.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.34'
$ws.Range("E41").Value = '  -0.70%  '

$ws.Range("E42").Value = '  +3.75%  '

$ws.Range("E43").Value = '  -0.60%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.92'
$ws.Range("E45").Value = '  +1.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.74'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.18'
$ws.Range("E47").Value = '  +1.19%  '

$ws.Range("E48").Value = '  -0.93%  '

$ws.Range("E49").Value = '  +0.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.616'
$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.566'
$ws.Range("E51").Value = '  +0.65%  '
